$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (2008年 and 2009年 data), shifting the remaining
# rows (2010年, 2011年) up so they become rows 2 and 3.
$ws.Range("A2:A3").EntireRow.Delete()
